$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 13 with the new run entry ("_ta_rq")
$ws.Range("A13").Value = "_ta_rq"

# Copy the date cell's existing style (numFmtId 16, "d-mmm") from the row above
# so a duplicate style entry isn't created in styles.xml
$ws.Range("B12").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B13").Value = 45909

$ws.Range("C13").Value = "Same as above, but with RQ kernel"
$ws.Range("D13").Value = "[-3.5, -4, -4.5, -5]"
$ws.Range("E13").Value = "[7, 10, 10, 10]"
$ws.Range("F13").Value = "<- rq kernel"

# Update the active selection to F13, matching the saved view state
$ws.Range("F13").Select()
